$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: new "Output" column ---
$ws.Range("X1").Copy()
$ws.Range("Y1").PasteSpecial(-4122)
$ws.Range("Y1").Value = "Output"

# --- Row 2: updated metric values ---
$ws.Range("Q2").Value = 0.9
$ws.Range("R2").Value = 0.1255144032921811
$ws.Range("S2").Value = 0.02073732718894009
$ws.Range("T2").Value = 0.9838709677419355
$ws.Range("U2").Value = 9
$ws.Range("V2").Value = 61
$ws.Range("W2").Value = 425
$ws.Range("X2").Value = 1

# --- Row 2: Output column payload ---
$ws.Range("Y2").Value = '[[''36823236'', ''35841079'', ''21608000'', ''25024623'', ''21900788''], [''35309311'', ''32071558'', ''28620221'', ''34843069'', ''37328872''], [''35858897'', ''36183106'', ''36698192'', ''35395830'', ''35278227''], [''35099348'', ''37918086''], [''29104502'', ''35505419'', ''37539053'', ''34882906'', ''37337578''], [''37498509'', ''36812711'', ''36998425'', ''36049782'', ''26420953''], [''37864199'', ''36983624''], [''27601441'', ''29303079'', ''32077830''], [''36544529'', ''37283714'', ''27176654'', ''36423545'', ''35902883''], [''36462157'', ''38164852'', ''25320545''], [''26691057'', ''30854976'', ''37282819''], [''26934883'', ''37329808'', ''37494836''], [''37798761'', ''28627746'', ''37114594'', ''35100287'', ''37852157''], [''24036368'', ''32827326'', ''33317249'', ''37626924'', ''27539153''], [''32928296'', ''37604091'', ''36029414'', ''37079286'', ''35563807''], [''38442726'', ''35932064'', ''27819130'', ''37271943'', ''32775366''], ["36093813'' , ''37365036", ''24415865''], [''31281379'', ''31591328''], [''37874430'', ''35484432'', ''32621283'', ''32267915''], [''27481213'', ''27495745'', ''22228987''], [''33554346'', ''34565478'', ''37167639''], [''30394698'', ''35659360'', ''30885249''], [''28297787'', ''35440002'', ''38013053'', ''27153445'', ''36552767''], [''31520775'', ''19642968'', ''19941447'', ''24157970'', ''20220785''], [''30106103'', ''30019838'', ''30248223'', ''27591095'', "27468195'' ,"], [''32880852'', ''31755595'', ''32962510''], [''38419761'', ''24457953'', ''31435712'', ''35883205'', ''36330597''], [''37345718'', ''32633558'', ''34186042'', ''37226279'', ''37742014''], [''31048740'', ''35084178'', ''28805262'', ''37084952''], [''25450712'', ''33143043'', ''31781620'', ''29497038'', ''32865681''], [''30375438'', ''35578365'', ''31592593''], [''36703021'', ''37298538'', ''37507381'', ''38105754'', ''38198960''], [''35199612'', ''37501214'', ''38159943''], [''37980535'', ''30635047'', ''32816183'', ''33935044'', ''35006577''], [''30223894'', ''23531302'', ''21200380'', ''22678013'', ''27500218''], [''37787035'', ''35036058'', ''30834983'', ''34518510'', ''33384547''], [''25613158'', ''32506351'', ''35293541'', ''26811102''], [''33293768'', ''34046793'', ''33145022'', ''31760709'', ''34012193''], [''25006200'', ''31147180'', ''36649501'', ''25722117'', ''28583205''], [''30048056'', ''33389680'', ''22533821'', ''34809720'', ''30675167''], [''23686136'', ''26464632'', ''22753359'', ''31245923'', ''37226801''], [''36359733'', ''36735212'', ''36411771'', ''36183364'', ''27908750''], [''33555212'', ''33348053'', ''30850965'', ''33060633'', ''33785758''], [''35435145'', ''34775285'', ''32374186'', ''26310071''], [''29266258'', ''22416179'', ''29286155'', ''23002959'', ''36674693''], [''33708992'', ''33342075'', ''42774189''], [''35075209'', ''38402155'', ''32389821'', ''33617844'', ''34016164''], [''19780871'', ''26369870'', ''33349209'', ''27805290'', ''30538216''], [''20943307'', ''24507685'', ''36941658'', ''32011099'', ''38356602''], [''32245503'', ''29861737'', ''30901579'', ''35186473''], [''32535098'', ''30810620'', ''37031397'', ''36610716'', ''25315017''], [''26876222'', ''33629329'', ''22182291'', ''30635054'', ''25445161''], [''27365486'', ''28382720'', ''37899504'', ''31452434'', ''34859376''], [''28494546'', ''35831878'', ''37133409'', ''30402110'', ''34174353''], [''24488807'', ''35315848'', ''37064607'', ''30700543'', ''33855751''], [''17869217'', ''26157751'', ''28544786'', ''23809154'', ''19684854''], [''34343826'', ''31508820'', ''17543295'', ''27981544'', ''28115038''], [''33078915'', ''27876093'', ''29130552'', ''34158112'', ''31768225''], [''32546282'', ''33855079'', ''24835895'', ''19686763''], [''27430962'', ''21289474'', ''30114694'', ''29145157'', ''34689832''], [''18176050'', ''32289375'', ''22609227'', ''29362496'', ''29264872''], [''24104560'', ''29150810'', ''34514810'', ''31550244'', ''25490388''], [''32357508'', ''34627378'', ''20178458'', ''22735033'', ''34712665''], [''19398237'', ''29398486'', ''21966929'', ''34256837'', ''28301671''], [''29512750'', ''28217916'', ''27028843'', ''26353075'', ''34416908''], [''29666468'', ''21326862'', ''25678175''], [''22927965'', ''33246509'', ''20367252'', ''23402186'', ''26861623''], [''20096523'', ''16469408'', ''32650827'', ''29290652''], [''19422480'', ''27280796'', ''29804359'', ''32801431'', ''31525863''], [''32867602'', ''25426285'', ''29223366'', ''38016185'', ''18486597''], [''28413861'', ''32603821'', ''16150456'', ''31552493'', ''30979485''], [''27960154'', ''24481983'', ''22320928'', ''31551525'', ''27543493''], [''28610623'', ''22167104'', ''31952158'', ''23855301'', ''27956229''], [''16233854'', ''23807535'', ''28585911'', ''38540717'', ''38413394''], [''24886681'', ''21356258'', ''22432472'', ''26125913'', ''32503644''], [''20025519'', ''29187911'', ''21537835'', ''25501583'', ''37964155''], [''27130910'', ''26389810'', ''19559033'', ''20824962'', ''26877613''], [''29110545'', ''30230766'', ''25473446'', ''25639333'', ''31982488''], [''14568247'', ''29293276'', ''23150236'', ''23782511'', ''22941276''], [''26153152'', ''18522880'', ''20663626'', ''25842730'', ''21810549''], [''25145391'', ''26643997'', ''38360740'', ''21770423'', ''22056375''], [''19399744'', ''37997442'', ''38474368'', ''28782844'', ''38160190''], [''34865351'', ''23638052'', ''19455046'', ''19056764'', ''24144538''], [''27938474'', ''19242797'', ''30453717'', ''19476212'', ''29074333''], [''25063290'', ''24115681'', ''38429940'', ''28035356'', ''24682874''], [''21239501'', ''33238629'', ''24967970'', ''17903050''], [''26194009'', ''24942592'', ''26839564'', ''28625529''], [''21145826'', ''20184502'', ''25268050'', ''26358689''], [''23849977'', ''20182823'', ''30631870'', ''21989829''], [''24815159'', ''26155209'', ''26599535'', ''25446803'', ''21269936''], [''24116127'', ''27878730'', ''21391433'', ''17903371'', ''17952094''], [''26983390'', ''24417995'', ''25890163'', ''28073411'', ''25654284''], [''25203804'', ''25518471'', ''15850515'', ''16108363'', ''38474048''], [''26540984''], [''23522251'', ''15650781'', ''24748190''], [''21526984'', ''30461638'', ''27339398''], [''29285145'', ''23980370'', ''31885614'', ''27405780'', ''34715410''], [''35721175'', ''34504910'', ''27785340'', ''34711695''], [''21516840'', ''37576406'', ''29844633''], [''34796180'']]'

# --- Rows 3-5: blank "Output" cells (text-typed empty, matching the rest of the row) ---
$ws.Range("Y3").Value = "'"
$ws.Range("Y3").Style = "Normal"
$ws.Range("Y4").Value = "'"
$ws.Range("Y4").Style = "Normal"
$ws.Range("Y5").Value = "'"
$ws.Range("Y5").Style = "Normal"
